# hmw 3 Postman task 1
# - A10 becomes the text value "1,23" (stored as shared string) instead of
#   the numeric value 1.23
# - the sheet's active selection moves to D21

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text interpretation so "1,23" is stored as a string (not parsed as
# a locale-specific number) while keeping the existing cell style.
$ws.Range("A10").NumberFormat = "@"
$ws.Range("A10").Value = "1,23"

# Move the selection to D21 (single cell, no prior multi-range selection).
$ws.Range("D21").Select()
